$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing column C values for rows 4-9 (id_energy_carrier precalculated) ---
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 12
$ws.Range("C8").Value = 13
$ws.Range("C9").Value = 14

# --- Add new rows 10-13, same shape as row 9 (id_scenario=1, id_region=9, unit=1, years=0.19) ---
$newCValues = @(15, 19, 25, 26)
for ($i = 0; $i -lt $newCValues.Length; $i++) {
    $r = 10 + $i
    $ws.Range("A$r").Value = 1
    $ws.Range("B$r").Value = 9
    $ws.Range("C$r").Value = $newCValues[$i]
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r`:AJ$r").Value = 0.19
    $ws.Range("E$r`:AJ$r").NumberFormat = $ws.Range("E9").NumberFormat
}

# --- Resize the table/list object to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:AJ13"))

# --- Update the active selection shown in the sheet view ---
$ws.Range("AG18").Select()
